$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: repurpose the old 172.16.0.0 network objects into host objects
$ws.Range("A2").Value = "Host_12.0.0.1"
$ws.Range("B2").Value = "12.0.0.1"

$ws.Range("A3").Value = "Host_12.0.0.3"
$ws.Range("B3").Value = "12.0.0.3"

# Add a new network object row, including group assignment
$ws.Range("A4").Value = "Net_15.0.0.0"
$ws.Range("B4").Value = "15.0.0.0/24"
$ws.Range("C4").Value = "Remote_User"

# Move the active selection to reflect where the user left off editing
$ws.Range("B10").Select()
